$d = $word.ActiveDocument

# --------------------------------------------------------------------
# 1. Split the "Short Description" text so the cell ends with a bold
#    colon: "The talk focused on generating evidence for operational
#    excellence, insights" -> "The talk focused on:"
# --------------------------------------------------------------------
$d.Content.Find.Execute(
    "The talk focused on generating evidence for operational excellence, insights",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "The talk focused on:", 2) | Out-Null

# --------------------------------------------------------------------
# 2. Remove the row that used to continue that sentence
#    ("and improvements in clinical practice, epidemiological
#    evidence using:") - find it by its cell content and delete it.
# --------------------------------------------------------------------
$t = $d.Tables.Item(3)

for ($i = 1; $i -le $t.Rows.Count; $i++) {
    $cellText = $t.Cell($i, 2).Range.Text
    if ($cellText -like "*and improvements in clinical practice, epidemiological evidence*") {
        $t.Rows.Item($i).Delete()
        break
    }
}

# --------------------------------------------------------------------
# 3. Re-sync every remaining row's cell widths with the table's grid
#    (1890/7380 dxa -> 1867/7178 dxa), matching the post-deletion
#    AutoFit recalculation Word performs on the surviving rows.
# --------------------------------------------------------------------
$newCol1 = 1867 / 20.0
$newCol2 = 7178 / 20.0

for ($i = 1; $i -le $t.Rows.Count; $i++) {
    $t.Cell($i, 1).Width = $newCol1
    $t.Cell($i, 2).Width = $newCol2
}
